$d = $word.ActiveDocument

# Anchor on the last existing paragraph of text ("Some of my functions...")
# and the (already existing) blank paragraph that immediately follows it.
$anchor = $d.Content
$anchor.Find.Execute("Some of my functions are quite long, but no time to go back and split out to cleaner code.")
$anchorPara = $anchor.Paragraphs(1)
$blankAfterAnchorIndex = $anchorPara.Index + 1

# --- New paragraph 1: VS Code / root paths / GitHub Pages note -------------
$d.Paragraphs($blankAfterAnchorIndex).Range.InsertParagraphAfter()
$para1Index = $blankAfterAnchorIndex + 1

$p1 = $d.Paragraphs($para1Index).Range
$p1.InsertAfter("Discovered that my site would work fine through VS Code, but due to using root paths, images and links were not being loaded correctly when I ran the site by opening up the html file or via GitHub pages")

$p1 = $d.Paragraphs($para1Index).Range
$p1.Collapse(0)
$p1.InsertAfter(". Had to modify paths")

$p1 = $d.Paragraphs($para1Index).Range
$p1.Collapse(0)
$p1.InsertAfter(" and pass modifier to js file to correct for nested files calling same script.")

# --- Blank separator paragraph ---------------------------------------------
$d.Paragraphs($para1Index).Range.InsertParagraphAfter()
$blankIndex = $para1Index + 1

# --- New paragraph 2: browser blocking type="module" / json note -----------
$d.Paragraphs($blankIndex).Range.InsertParagraphAfter()
$para2Index = $blankIndex + 1

$p2 = $d.Paragraphs($para2Index).Range
$p2.InsertAfter("Discovered browser was blocking the js file of type=”module” (also blocks import statements). This worked fine via VS code, but found error when pushed to github pages. Ended up having to load json object to local storage by creating json object in js file instead of")

$p2 = $d.Paragraphs($para2Index).Range
$p2.Collapse(0)
$p2.InsertAfter(" importing from")

$p2 = $d.Paragraphs($para2Index).Range
$p2.Collapse(0)
$p2.InsertAfter(" json file. (no time to look into further).")

# --- Trailing blank paragraph (mirrors pressing Enter after the last
#     sentence), keeping the document ending on an empty paragraph as before.
$d.Paragraphs($para2Index).Range.InsertParagraphAfter()
